# feat: Sum after sale fee
#
# The "receivedAmount" column (D) for GHST-denominated Baazaar sales stored
# the gross sale price. This updates it to the net amount the seller
# actually received after the Baazaar marketplace fee is deducted.
#
# Sales from 2021-06-16 onward (rows 2-51) used a 3.5% fee (net multiplier
# 0.965); earlier sales (rows 55+, before that date) used a 3% fee (net
# multiplier 0.97). Rows whose receivedCurrency (column E) is not GHST
# (e.g. NFT transfers recorded with receivedAmount = 1) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($row = 2; $row -le $lastRow; $row++) {
    $currency = $ws.Cells.Item($row, 5).Value2   # column E: receivedCurrency
    if ($currency -eq "GHST") {
        $amount = $ws.Cells.Item($row, 4).Value2  # column D: receivedAmount
        if ($row -le 51) {
            $factor = 0.965
        } else {
            $factor = 0.97
        }
        $ws.Cells.Item($row, 4).Value2 = $amount * $factor
    }
}
